$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are being dropped from the costing table:
# old row 14 ("Operation costs") and old row 12 ("Flagging system").
# Delete bottom-up so the remaining row numbers stay predictable.
$ws.Range("A14:C14").EntireRow.Delete()
$ws.Range("A12:C12").EntireRow.Delete()

# Fix up the ID numbering for the two rows that shifted into the gap.
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12

# Updated feature costs.
$ws.Range("C2").Value = 300
$ws.Range("C3").Value = 300
$ws.Range("C4").Value = 300
$ws.Range("C7").Value = 300
$ws.Range("C8").Value = 600
$ws.Range("C10").Value = 2000
$ws.Range("C11").Value = 300
$ws.Range("C12").Value = 300

# Selection moves from D1 to C1.
$ws.Range("C1").Select()

# Page setup now targets A4 portrait printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
